$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("organizations")

# The table used to start at column B (STT, Tên đơn vị, Mã quản lý, Mã cấp trên).
# It now starts at column A, and two new columns (Tên viết tắt, Mô tả) are
# inserted after "Tên đơn vị". Net effect: drop the empty leading column A,
# then insert two fresh columns at C:D (they inherit the neighbouring format).
$ws.Columns("A:A").Delete()
$ws.Columns("C:D").Insert()

# New header cells for the two inserted columns.
$ws.Range("C3").Value = "Tên viết tắt"
$ws.Range("D3").Value = "Mô tả"

# Shorten the existing "Phòng Tự động hóa" name now that there is a
# dedicated abbreviation column.
$ws.Range("B6").Value = "Phòng Tự động"

# Fill in abbreviations for each existing department.
$ws.Range("C4").Value = "HC"
$ws.Range("C5").Value = "DV"
$ws.Range("C6").Value = "TĐ"

# Fill in descriptions for each existing department.
$ws.Range("D4").Value = "Mô tả phòng hành chính"
$ws.Range("D5").Value = "Mô tả phòng dịch vụ"
$ws.Range("D6").Value = "Mô tả phòng tự động"

# Append a brand-new department row.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Phòng Giải lao"
$ws.Range("C7").Value = "GL"
$ws.Range("D7").Value = "Mô tả giải lao"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 1

# Match the column widths as closely as the engine's rounding allows.
$ws.Columns("A").ColumnWidth = 6.5
$ws.Columns("B").ColumnWidth = 16
$ws.Columns("C").ColumnWidth = 8.83333333333333
$ws.Columns("D").ColumnWidth = 22.6666666666667
$ws.Columns("E").ColumnWidth = 9.33333333333333
$ws.Columns("F").ColumnWidth = 9.83333333333333

[void]$ws.Range("E4").Select()
